# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) and updates recomputed metric values (D, E, F)
# across the two per-patient blocks (rows 2-8 and 9-15) of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell H1 -------------------------------------------------
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# --- Updated numeric values (D, E, F) ------------------------------------
$ws.Range("D2").Value = 0.7834205856158361
$ws.Range("E2").Value = 0.7834205856158361

$ws.Range("D3").Value = 0.5466942929768367
$ws.Range("E3").Value = 0.5466942929768367

$ws.Range("D4").Value = 0.3278335971385047
$ws.Range("E4").Value = 0.6721664028614953

$ws.Range("D5").Value = 0.5744591766741292
$ws.Range("E5").Value = 0.4255408233258708

$ws.Range("D6").Value = 0.570834886719044
$ws.Range("E6").Value = 0.429165113280956

$ws.Range("F8").Value = 0.9565996527671814

# --- New "Label" values (H2:H15) -----------------------------------------
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
